$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.047.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.96%  "

$ws.Range("D3").Value = "'2.463.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.18%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "'574.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.91%  "

$ws.Range("D6").Value = "'145.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.82%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  +2.42%  "

$ws.Range("D9").Value = "'2.461.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.34%  "

$ws.Range("E10").Value = "  +6.12%  "

$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").Value = "'5.23"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.54%  "

$ws.Range("D13").Value = "'0.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.33%  "

$ws.Range("D14").Value = "'27.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.65%  "

$ws.Range("D15").Value = "'0.0000178"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +8.69%  "

$ws.Range("D16").Value = "'2.958.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.05%  "

$ws.Range("D17").Value = "'62.992.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.03%  "

$ws.Range("D18").Value = "'2.448.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.66%  "

$ws.Range("D19").Value = "'7.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.05%  "

$ws.Range("D20").Value = "'10.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.22%  "

$ws.Range("D21").Value = "'327.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.20%  "

$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").Value = "'2.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.47%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.26%  "

$ws.Range("D25").Value = "'65.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("D26").Value = "'623.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.08%  "

$ws.Range("E27").Value = "  +10.28%  "

$ws.Range("D28").Value = "'8.39"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.86%  "

$ws.Range("D29").Value = "'0.0₃0982"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.73%  "

$ws.Range("D31").Value = "'8.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.61%  "

$ws.Range("D32").Value = "'1.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.97%  "

$ws.Range("D33").Value = "'0.137"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.32%  "

$ws.Range("D34").Value = "'1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.33%  "

$ws.Range("D35").Value = "'1.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.40%  "

$ws.Range("D36").Value = "'0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("E37").Value = "  +5.44%  "

$ws.Range("D38").Value = "'0.373"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.36%  "

$ws.Range("D39").Value = "'152.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("D40").Value = "'5.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.23%  "

$ws.Range("D41").Value = "'18.66"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.53%  "

$ws.Range("E42").Value = "  +15.94%  "

$ws.Range("E43").Value = "  +7.38%  "

$ws.Range("D44").Value = "'42.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.77%  "

$ws.Range("E45").Value = "  -0.03%  "

$ws.Range("D46").Value = "'0.0₆0286"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.81%  "

$ws.Range("D47").Value = "'144.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.03%  "

$ws.Range("D48").Value = "'3.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.71%  "

$ws.Range("D49").Value = "'20.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.73%  "

$ws.Range("D50").Value = "'0.602"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.13%  "

$ws.Range("D51").Value = "'0.0517"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.58%  "
